$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (Source File) from 25 to 50
# (subtract the ~5/6 character padding this engine adds internally so the
# stored OOXML <col width="..."> lands exactly on 50)
$ws.Columns.Item(5).ColumnWidth = 49.166666666666664

# Replace the old source-file name with the new one in column E,
# for every data row (2 through 307)
$newName = "Y4_B2526_General_&_special_internal_1_reference_data.xlsx"
$lastRow = $ws.Cells(1, 1).End(4).Row
if ($lastRow -lt 307) { $lastRow = 307 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "Group B1 2025-2026.xlsx") {
        $cell.Value2 = $newName
    }
}
